$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) values for new columns O and P, using same style as existing header cells (N1 etc.)
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)

$ws.Range("O1").Value = 14
$ws.Range("P1").Value = 15

# Data rows 2-6 for new columns O and P
$ws.Range("O2").Value = -1.281641432976643
$ws.Range("P2").Value = -1.12762283105276

$ws.Range("O3").Value = -0.4560821420107516
$ws.Range("P3").Value = -0.4109627019024589

$ws.Range("O4").Value = 0.06196656500936742
$ws.Range("P4").Value = 0.03749380507472993

$ws.Range("O5").Value = 0.4341605616485296
$ws.Range("P5").Value = 0.4135709484694772

$ws.Range("O6").Value = -0.06451975945486933
$ws.Range("P6").Value = -0.05826341342551093
